$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text in D/E columns from Excel auto-number conversion
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '97.999.62'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '3.372.39'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '253.28'
$ws.Range("E5").Value = '  -1.99%  '
$ws.Range("D6").Value = '661.36'
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").Value = '1.41'
$ws.Range("E7").Value = '  -4.05%  '
$ws.Range("D8").Value = '0.425'
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '1.04'
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("D11").Value = '3.366.27'
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("E12").Value = '  -3.28%  '
$ws.Range("D13").Value = '41.66'
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("D14").Value = '97.553.13'
$ws.Range("D15").Value = '6.09'
$ws.Range("E15").Value = '  -5.50%  '
$ws.Range("E16").Value = '  -5.53%  '
$ws.Range("D17").Value = '4.007.55'
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = '8.80'
$ws.Range("E18").Value = '  -6.29%  '
$ws.Range("D19").Value = '3.376.62'
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("D20").Value = '18.01'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '0.523'
$ws.Range("E21").Value = '  -14.63%  '
$ws.Range("D22").Value = '10.95'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").Value = '511.96'
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("E24").Value = '  -4.39%  '
$ws.Range("D25").Value = '7.03'
$ws.Range("E25").Value = '  +8.61%  '
$ws.Range("D26").Value = '0.0000200'
$ws.Range("E26").Value = '  -3.89%  '
$ws.Range("D27").Value = '96.79'
$ws.Range("E27").Value = '  -5.09%  '
$ws.Range("D28").Value = '12.35'
$ws.Range("E28").Value = '  -7.04%  '
$ws.Range("D29").Value = '11.37'
$ws.Range("E29").Value = '  -5.56%  '
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '0.142'
$ws.Range("E31").Value = '  -8.25%  '
$ws.Range("D32").Value = '0.186'
$ws.Range("E32").Value = '  -6.45%  '
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  +7.46%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").Value = '0.561'
$ws.Range("E35").Value = '  -5.01%  '
$ws.Range("D36").Value = '28.69'
$ws.Range("E36").Value = '  -5.34%  '
$ws.Range("D37").Value = '8.06'
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D38").Value = '1.52'
$ws.Range("E38").Value = '  +2.58%  '
$ws.Range("D39").Value = '527.77'
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("D40").Value = '0.152'
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").Value = '24.41'
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '2.35'
$ws.Range("E43").Value = '  +9.71%  '
$ws.Range("D44").Value = '0.855'
$ws.Range("E44").Value = '  -5.03%  '
$ws.Range("D45").Value = '1.74'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").Value = '0.0427'
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("D47").Value = '3.66'
$ws.Range("E47").Value = '  -2.42%  '
$ws.Range("D48").Value = '5.63'
$ws.Range("E48").Value = '  -7.53%  '
$ws.Range("D49").Value = '8.56'
$ws.Range("E49").Value = '  -9.21%  '
$ws.Range("D50").Value = '55.06'
$ws.Range("E50").Value = '  +6.82%  '
$ws.Range("D51").Value = '3.18'
$ws.Range("E51").Value = '  -5.67%  '

# Restore original (default) style so no stray format attributes are left on cells
$dataRange.Style = "Normal"
